$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 204: this pushes the existing rows 204-222
# down to 205-223 (and the sheet dimension grows to R223), matching the
# weekly update where a brand-new record is prepended to the data block.
$ws.Rows("204:204").Insert()

# Populate the newly inserted row 204 with this week's record.
$ws.Range("A204").Value = 7
$ws.Range("B204").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C204").Value = "Ñuble"
$ws.Range("D204").Value = 45223
$ws.Range("E204").Value = 16
$ws.Range("F204").Value = 100112037
$ws.Range("G204").Value = "Cebollín"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 150
$ws.Range("K204").Value = 5000
$ws.Range("L204").Value = 5000
$ws.Range("M204").Value = 5000
$ws.Range("N204").Value = "$/paquete 36 unidades"
$ws.Range("O204").Value = "Provincia de Diguillín"
$ws.Range("P204").Value = 139
$ws.Range("Q204").Value = 36
$ws.Range("R204").Value = "Hortaliza"
